# The edit re-sorts the weekly price records (rows 2-41) into a new
# (date-ascending-ish) order. Every row keeps the same "constant" columns
# (A,B,C,E,F,G,H,I,J,K,Q,T) and only the per-record columns
# (D,L,M,N,O,P,R,S) are permuted between rows. Concretely, the new content
# of row <r> equals the OLD content (pre-edit) of row <map[r]>.
#
# map[new_row] = old_row_that_supplies_the_data
$map = @{
    2=11;  3=5;   4=30;  5=20;  6=25;  7=29;  8=37;  9=12;  10=14;
    11=24; 12=31; 13=2;  14=40; 15=17; 16=7;  17=38; 18=16; 19=4;
    20=26; 21=39; 22=32; 23=27; 24=19; 25=15; 26=22; 27=3;  28=13;
    29=36; 30=10; 31=21; 32=9;  33=35; 34=23; 35=6;  36=34; 37=8;
    38=33; 39=18; 40=41; 41=28
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the per-record data which gets permuted.
$cols = @(4, 12, 13, 14, 15, 16, 18, 19)   # D, L, M, N, O, P, R, S

# 1) Snapshot every relevant cell's current (pre-edit) value, keyed by row.
$snapshot = @{}
for ($r = 2; $r -le 41; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write back according to the permutation map. Because we already have
#    a full snapshot, write order doesn't matter.
foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}
